$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.029774912949598
$ws.Range("D2").Value = 1.032419332480934
$ws.Range("E2").Value = 1.038825699602864
$ws.Range("F2").Value = 1.048033410507566
$ws.Range("I2").Value = 1.033087251084681
$ws.Range("J2").Value = 1.034919824371565
$ws.Range("K2").Value = 1.035224617337204
$ws.Range("L2").Value = 1.041612629256448
$ws.Range("M2").Value = 1.050794385074352
$ws.Range("N2").Value = 1.015627508519687

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.030663321431301
$ws.Range("D3").Value = 1.033221146077235
$ws.Range("E3").Value = 1.039622013149809
$ws.Range("F3").Value = 1.04892939061912
$ws.Range("I3").Value = 1.033193230826109
$ws.Range("J3").Value = 1.035449719230273
$ws.Range("K3").Value = 1.035835217861591
$ws.Range("L3").Value = 1.042219072121332
$ws.Range("M3").Value = 1.051502111559733
$ws.Range("N3").Value = 1.015804924840287

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031238729547898
$ws.Range("D4").Value = 1.033740786047126
$ws.Range("E4").Value = 1.040138149967664
$ws.Range("F4").Value = 1.049510103153897
$ws.Range("I4").Value = 1.03326049611559
$ws.Range("J4").Value = 1.035792513870385
$ws.Range("K4").Value = 1.0362304775236
$ws.Range("L4").Value = 1.042611694556565
$ws.Range("M4").Value = 1.051960375538809
$ws.Range("N4").Value = 1.015919644888917

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.031480761043106
$ws.Range("D5").Value = 1.033959435661628
$ws.Range("E5").Value = 1.040355340196899
$ws.Range("F5").Value = 1.049754461231469
$ws.Range("I5").Value = 1.033288460060124
$ws.Range("J5").Value = 1.035936603739297
$ws.Range("K5").Value = 1.036396681729117
$ws.Range("L5").Value = 1.042776802783341
$ws.Range("M5").Value = 1.052153104253435
$ws.Range("N5").Value = 1.015967853557257

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.031521406775124
$ws.Range("D6").Value = 1.033996159161883
$ws.Range("E6").Value = 1.040391819442016
$ws.Range("F6").Value = 1.049795503261074
$ws.Range("I6").Value = 1.033293136876411
$ws.Range("J6").Value = 1.035960795796598
$ws.Range("K6").Value = 1.036424590280124
$ws.Range("L6").Value = 1.04280452806037
$ws.Range("M6").Value = 1.052185468579982
$ws.Range("N6").Value = 1.015975946840116

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.031241963075784
$ws.Range("D7").Value = 1.033743706897807
$ws.Range("E7").Value = 1.04014105126556
$ws.Range("F7").Value = 1.049513367390573
$ws.Range("I7").Value = 1.033260871007246
$ws.Range("J7").Value = 1.035794439290257
$ws.Range("K7").Value = 1.036232698208443
$ws.Range("L7").Value = 1.042613900546928
$ws.Range("M7").Value = 1.051962950498327
$ws.Range("N7").Value = 1.015920289133235

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030075040375917
$ws.Range("D8").Value = 1.032690139818247
$ws.Range("E8").Value = 1.039094636626097
$ws.Range("F8").Value = 1.048336012708718
$ws.Range("I8").Value = 1.033123338459763
$ws.Range("J8").Value = 1.035098921245864
$ws.Range("K8").Value = 1.03543093865333
$ws.Range("L8").Value = 1.041817534491977
$ws.Range("M8").Value = 1.051033497886164
$ws.Range("N8").Value = 1.015687483528324

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.028023034378818
$ws.Range("D9").Value = 1.030839919320245
$ws.Range("E9").Value = 1.037257446880323
$ws.Range("F9").Value = 1.046268742833465
$ws.Range("I9").Value = 1.032870976806158
$ws.Range("J9").Value = 1.033872747981782
$ws.Range("K9").Value = 1.034019422339923
$ws.Range("L9").Value = 1.040415928153544
$ws.Range("M9").Value = 1.049398182545347
$ws.Range("N9").Value = 1.015276655972096

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026657976248269
$ws.Range("D10").Value = 1.029610772489593
$ws.Range("E10").Value = 1.036037274880895
$ws.Range("F10").Value = 1.044895627857534
$ws.Range("I10").Value = 1.032696040116329
$ws.Range("J10").Value = 1.033054978318899
$ws.Range("K10").Value = 1.033079355423967
$ws.Range("L10").Value = 1.039482740281818
$ws.Range("M10").Value = 1.048309741603813
$ws.Range("N10").Value = 1.015002397815572

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.026067606609266
$ws.Range("D11").Value = 1.029079586828132
$ws.Range("E11").Value = 1.035510043749927
$ws.Range("F11").Value = 1.044302276960854
$ws.Range("I11").Value = 1.03261870997829
$ws.Range("J11").Value = 1.032700813648829
$ws.Range("K11").Value = 1.032672536542329
$ws.Range("L11").Value = 1.039078964605488
$ws.Range("M11").Value = 1.04783887302388
$ws.Range("N11").Value = 1.014883558062659

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025848424971805
$ws.Range("D12").Value = 1.028882439184785
$ws.Range("E12").Value = 1.03531437527732
$ws.Range("F12").Value = 1.04408206459534
$ws.Range("I12").Value = 1.032589749115616
$ws.Range("J12").Value = 1.032569252525049
$ws.Range("K12").Value = 1.032521462720734
$ws.Range("L12").Value = 1.038929030854407
$ws.Range("M12").Value = 1.047664038097418
$ws.Range("N12").Value = 1.01483940348551

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.025895435271455
$ws.Range("D13").Value = 1.028924720830649
$ws.Range("E13").Value = 1.035356339170591
$ws.Range("F13").Value = 1.044129292514795
$ws.Range("I13").Value = 1.03259597203467
$ws.Range("J13").Value = 1.032597473202886
$ws.Range("K13").Value = 1.032553866880649
$ws.Range("L13").Value = 1.038961190027011
$ws.Range("M13").Value = 1.047701537743474
$ws.Range("N13").Value = 1.014848875337237

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.026049486777224
$ws.Range("D14").Value = 1.029063287302293
$ws.Range("E14").Value = 1.035493866282993
$ws.Range("F14").Value = 1.044284070362635
$ws.Range("I14").Value = 1.032616320895841
$ws.Range("J14").Value = 1.032689938933882
$ws.Range("K14").Value = 1.032660047974239
$ws.Range("L14").Value = 1.039066570078379
$ws.Range("M14").Value = 1.047824419746439
$ws.Range("N14").Value = 1.014879908477887

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.026144417336325
$ws.Range("D15").Value = 1.02914868371292
$ws.Range("E15").Value = 1.035578623672863
$ws.Range("F15").Value = 1.044379458613344
$ws.Range("I15").Value = 1.032628827112561
$ws.Range("J15").Value = 1.032746909018462
$ws.Range("K15").Value = 1.032725474569643
$ws.Range("L15").Value = 1.039131504408855
$ws.Range("M15").Value = 1.047900140277919
$ws.Range("N15").Value = 1.014899027416474

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026697172199692
$ws.Range("D16").Value = 1.029646047649559
$ws.Range("E16").Value = 1.036072289033044
$ws.Range("F16").Value = 1.044935032390928
$ws.Range("I16").Value = 1.032701138994243
$ws.Range("J16").Value = 1.033078481791891
$ws.Range("K16").Value = 1.033106359744803
$ws.Range("L16").Value = 1.039509543994204
$ws.Range("M16").Value = 1.048341000862409
$ws.Range("N16").Value = 1.015010283083196

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.027044091625089
$ws.Range("D17").Value = 1.029958311326356
$ws.Range("E17").Value = 1.036382251002301
$ws.Range("F17").Value = 1.045283856186034
$ws.Range("I17").Value = 1.032746075351979
$ws.Range("J17").Value = 1.033286452054123
$ws.Range("K17").Value = 1.033345343150073
$ws.Range("L17").Value = 1.039746759708821
$ws.Range("M17").Value = 1.048617658184347
$ws.Range("N17").Value = 1.015080048603846

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.027246512227213
$ws.Range("D18").Value = 1.030140549971705
$ws.Range("E18").Value = 1.036563153696551
$ws.Range("F18").Value = 1.045487436505051
$ws.Range("I18").Value = 1.032772133348316
$ws.Range("J18").Value = 1.033407751238477
$ws.Range("K18").Value = 1.033484760837482
$ws.Range("L18").Value = 1.039885152559215
$ws.Range("M18").Value = 1.048779069448559
$ws.Range("N18").Value = 1.015120733501683

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.027315543977904
$ws.Range("D19").Value = 1.030202705640728
$ws.Range("E19").Value = 1.036624854959247
$ws.Range("F19").Value = 1.045556871976095
$ws.Range("I19").Value = 1.032780992546867
$ws.Range("J19").Value = 1.033449110012835
$ws.Range("K19").Value = 1.03353230247782
$ws.Range("L19").Value = 1.039932345813039
$ws.Range("M19").Value = 1.048834113550859
$ws.Range("N19").Value = 1.015134604604596

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.027006863372419
$ws.Range("D20").Value = 1.029924797979233
$ws.Range("E20").Value = 1.036348983924289
$ws.Range("F20").Value = 1.04524641854847
$ws.Range("I20").Value = 1.032741269888166
$ws.Range("J20").Value = 1.033264139461652
$ws.Range("K20").Value = 1.033319700133797
$ws.Range("L20").Value = 1.039721305704981
$ws.Range("M20").Value = 1.048587971144924
$ws.Range("N20").Value = 1.01507256426353

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.026004119458217
$ws.Range("D21").Value = 1.029022478545173
$ws.Range("E21").Value = 1.035453363307348
$ws.Range("F21").Value = 1.044238487034963
$ws.Range("I21").Value = 1.032610335201274
$ws.Range("J21").Value = 1.032662710308642
$ws.Range("K21").Value = 1.032628779261499
$ws.Range("L21").Value = 1.039035536985211
$ws.Range("M21").Value = 1.047788232207997
$ws.Range("N21").Value = 1.014870770324716

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025374279688353
$ws.Range("D22").Value = 1.028456071560941
$ws.Range("E22").Value = 1.034891227628996
$ws.Range("F22").Value = 1.043605829253905
$ws.Range("I22").Value = 1.032526640409728
$ws.Range("J22").Value = 1.032284518774835
$ws.Range("K22").Value = 1.032194583598072
$ws.Range("L22").Value = 1.03860463673616
$ws.Range("M22").Value = 1.047285790836067
$ws.Range("N22").Value = 1.014743824090581

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.025708109966431
$ws.Range("D23").Value = 1.028756247139688
$ws.Range("E23").Value = 1.03518913326497
$ws.Range("F23").Value = 1.04394111122482
$ws.Range("I23").Value = 1.032571138359022
$ws.Range("J23").Value = 1.032485009530749
$ws.Range("K23").Value = 1.032424738235042
$ws.Range("L23").Value = 1.038833039100891
$ws.Range("M23").Value = 1.047552107402952
$ws.Range("N23").Value = 1.014811127214245

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.027023685006206
$ws.Range("D24").Value = 1.02993994090782
$ws.Range("E24").Value = 1.036364015554087
$ws.Range("F24").Value = 1.045263334643161
$ws.Range("I24").Value = 1.032743441742052
$ws.Range("J24").Value = 1.033274221581509
$ws.Range("K24").Value = 1.033331287039039
$ws.Range("L24").Value = 1.039732807184399
$ws.Range("M24").Value = 1.048601385311409
$ws.Range("N24").Value = 1.015075946139992

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028553013836689
$ws.Range("D25").Value = 1.031317489048049
$ws.Range("E25").Value = 1.037731597173316
$ws.Range("F25").Value = 1.046802296489837
$ws.Range("I25").Value = 1.032937401105401
$ws.Range("J25").Value = 1.03418980485091
$ws.Range("K25").Value = 1.034384172547119
$ws.Range("L25").Value = 1.040778068344
$ws.Range("M25").Value = 1.049820645440163
$ws.Range("N25").Value = 1.015382932445168
